$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.292685389518738
$ws.Range("B1").Value = 1.445991158485413
$ws.Range("C1").Value = 1.313401341438293
$ws.Range("D1").Value = 1.443878293037415
$ws.Range("E1").Value = 1.122377038002014
